$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 16; everything from row 16 downward shifts down by one.
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new record.
$ws.Cells.Item(16, 1).Value = 3
$ws.Cells.Item(16, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(16, 3).Value = 'Coquimbo'
$ws.Cells.Item(16, 4).Value = 44550
$ws.Cells.Item(16, 5).Value = 5
$ws.Cells.Item(16, 6).Value = 100112026
$ws.Cells.Item(16, 7).Value = 'Haba'
$ws.Cells.Item(16, 8).Value = 'Sin especificar'
$ws.Cells.Item(16, 9).Value = 'Primera'
$ws.Cells.Item(16, 10).Value = 75
$ws.Cells.Item(16, 11).Value = 7500
$ws.Cells.Item(16, 12).Value = 8000
$ws.Cells.Item(16, 13).Value = 7767
$ws.Cells.Item(16, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(16, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(16, 16).Value = 311
$ws.Cells.Item(16, 17).Value = 25
$ws.Cells.Item(16, 18).Value = 'Hortaliza'

# Match the date cell style used by the rest of column D.
$ws.Cells.Item(16, 4).NumberFormat = $ws.Cells.Item(17, 4).NumberFormat
